$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.028.10'
$ws.Range('E2').Value = '  -1.61%  '
$ws.Range('D3').Value = '1.628.62'
$ws.Range('E3').Value = '  -1.59%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9963'
$ws.Range('E5').Value = '  -0.54%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9970'
$ws.Range('E6').Value = '  -0.42%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3938'
$ws.Range('E7').Value = '  +0.51%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3865'
$ws.Range('E8').Value = '  -1.13%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '50.41'
$ws.Range('E9').Value = '  +0.65%  '
$ws.Range('E10').Value = '  -0.38%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.372'
$ws.Range('E11').Value = '  -0.88%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08504'
$ws.Range('E12').Value = '  -0.63%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '24.06'
$ws.Range('E13').Value = '  -3.31%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.127'
$ws.Range('E14').Value = '  -1.32%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.660'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.00001292'
$ws.Range('E16').Value = '  -1.06%  '
$ws.Range('D17').Value = '1.620.87'
$ws.Range('E17').Value = '  -2.34%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '94.21'
$ws.Range('E18').Value = '  +1.11%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06932'
$ws.Range('E19').Value = '  -0.29%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '20.24'
$ws.Range('E20').Value = '  -3.19%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.895'
$ws.Range('E21').Value = '  -1.67%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9966'
$ws.Range('E22').Value = '  -0.47%  '
$ws.Range('E23').Value = '  -2.27%  '
$ws.Range('D24').Value = '24.029.45'
$ws.Range('E24').Value = '  -1.56%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.463'
$ws.Range('E25').Value = '  +5.32%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.884'
$ws.Range('E26').Value = '  +3.46%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.34'
$ws.Range('E27').Value = '  -1.59%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '157.30'
$ws.Range('E28').Value = '  -0.99%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '140.94'
$ws.Range('E29').Value = '  -3.14%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.317'
$ws.Range('E30').Value = '  -7.10%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.088'
$ws.Range('E31').Value = '  -1.92%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.494'
$ws.Range('E32').Value = '  -2.61%  '
$ws.Range('D33').Value = '1.804.77'
$ws.Range('E33').Value = '  -1.85%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08206'
$ws.Range('E34').Value = '  +1.12%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.003'
$ws.Range('E35').Value = '  -1.04%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02928'
$ws.Range('E36').Value = '  -3.00%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.697'
$ws.Range('E37').Value = '  -2.42%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2690'
$ws.Range('E38').Value = '  -2.61%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '10.53'
$ws.Range('E39').Value = '  +3.15%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.09186'
$ws.Range('E40').Value = '  -2.70%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '13.77'
$ws.Range('E41').Value = '  +2.73%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.7617'
$ws.Range('E42').Value = '  -2.29%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.428'
$ws.Range('E43').Value = '  -3.79%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '16.20'
$ws.Range('E44').Value = '  -0.59%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6979'
$ws.Range('E45').Value = '  -0.75%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.490'
$ws.Range('E46').Value = '  -2.56%  '
$ws.Range('E47').Value = '  -1.45%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.9967'
$ws.Range('E48').Value = '  -0.37%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.08318'
$ws.Range('E49').Value = '  -3.19%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '136.85'
$ws.Range('E50').Value = '  +0.26%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.222'
$ws.Range('E51').Value = '  -6.35%  '
